$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 510, shifting existing rows 510-529 down to 511-530.
$ws.Rows.Item(510).Insert()

# Populate the new row with the Voteview dataset entry.
$ws.Range("A510").Value = "Voteview: Congressional Roll-Call Votes Database"
$ws.Range("B510").Value = "parties and politicians"
$ws.Range("C510").Value = "https://voteview.com/data"
$ws.Range("D510").Value = "roll-call voting, member ideology"
$ws.Range("E510").Value = "US"
$ws.Range("F510").Value = 0
$ws.Range("G510").Value = 0
$ws.Range("H510").Value = 0
$ws.Range("I510").Value = 0
$ws.Range("J510").Value = 1
$ws.Range("K510").Value = 1721
$ws.Range("L510").Value = 2022
$ws.Range("M510").Value = "online"
$ws.Range("N510").Value = "no"
$ws.Range("O510").Value = 1
$ws.Range("AB510").Value = 20230301

# Add the hyperlink for the link cell, then restore the usual "link" column
# cell style (matching the rest of column C) since Hyperlinks.Add otherwise
# creates its own ad-hoc style.
$linkStyle = $ws.Range("C511").Style
$ws.Hyperlinks.Add($ws.Range("C510"), "https://voteview.com/data")
$ws.Range("C510").Style = $linkStyle

# Reflect the editor's final scroll/selection position.
$excel.ActiveWindow.ScrollRow = 527
$ws.Range("A535").Select()
